# edit.ps1
#
# Reproduces the slide1.xml changes from the commit "edited and inserted
# Geoff's backprop figures and captions":
#
#   1. The slide no longer shows the slide-master's shapes
#      (p:sld/@showMasterSp="0"  <->  Slide.DisplayMasterShapes = False).
#   2. The slide gets its own solid-colour background (near-white,
#      RGB 255,255,254 / hex FFFFFE) instead of inheriting the master's
#      background (p:cSld/p:bg/p:bgPr/a:solidFill/a:srgbClr val="FFFFFE").
#
# (The diff also renumbers the legacy VML `spid` attributes on the ten
# equation OLE objects, e.g. _x0000_s1134 -> _x0000_s1155. That id is
# internal VML/legacy-drawing bookkeeping that PowerPoint assigns itself
# when shapes are (re)created; it is not exposed anywhere on the
# Shape/OLEFormat object model, so it cannot be targeted from COM/VBA
# automation and is intentionally left untouched here.)

function Get-RGBValue($r, $g, $b) {
    # Mirrors VBA's RGB() helper: packs R,G,B into the single OLE_COLOR
    # integer (0x00BBGGRR) that ForeColor.RGB expects.
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Stop inheriting shapes placed on the slide master.
$s.DisplayMasterShapes = 0

# Give the slide its own background fill rather than the master's.
$s.FollowMasterBackground = 0
$s.Background.Fill.Solid()
$s.Background.Fill.ForeColor.RGB = Get-RGBValue 255 255 254
